$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$summary = $wb.Worksheets.Item("Summary")

# Fix typo in the "Rendering provider" example value (was "Provider, Rending [123]")
$summary.Range("C3").Value = "Provider, Rendering [123]"

# Lowercase the "Facilities" example value
$summary.Range("C4").Value = "facilities"

# Add clarifying notes in column F
$summary.Range("F2").Value = "Items must be written exactly the same"
$summary.Range("F4").Value = "Payer and templates must match the tab name"

# Make Summary the active sheet/tab and leave the selection where the
# user's last edit ended up (after typing into F4 and pressing Enter).
$summary.Activate()
$summary.Range("F5").Select()
